$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column AP (42) from width 23 to width 20
$ws.Columns.Item(42).ColumnWidth = 19.1667

# Remove stray spaces from HD identifiers in column AP
$ws.Range("AP2").Value = "HD196761"
$ws.Range("AP3").Value = "HD172051"
$ws.Range("AP4").Value = "HD130992"
$ws.Range("AP5").Value = "HD220339"
$ws.Range("AP6").Value = "HD10700"
$ws.Range("AP7").Value = "HD35854"
$ws.Range("AP8").Value = "HD82342"
$ws.Range("AP9").Value = "HD52919"
$ws.Range("AP10").Value = "HD189567"
$ws.Range("AP11").Value = "HD156274"
$ws.Range("AP12").Value = "HD113538"
$ws.Range("AP13").Value = "HD142709"
$ws.Range("AP14").Value = "HD136352"
$ws.Range("AP15").Value = "HD115617"
$ws.Range("AP16").Value = "HD102438"
$ws.Range("AP17").Value = "HD55"
$ws.Range("AP18").Value = "HD170493"
$ws.Range("AP19").Value = "HD196877"
$ws.Range("AP20").Value = "HD101581"
$ws.Range("AP21").Value = "HD192961"
$ws.Range("AP22").Value = "HD21209A,HD21209"
$ws.Range("AP23").Value = "HD56533"
$ws.Range("AP24").Value = "HD22496"
$ws.Range("AP25").Value = "HD2025"
$ws.Range("AP26").Value = "HD20794"
$ws.Range("AP27").Value = "HD26965"
$ws.Range("AP28").Value = "HD21749"
$ws.Range("AP29").Value = "HD160346"
$ws.Range("AP30").Value = "HD157347"
$ws.Range("AP31").Value = "HD154363"
$ws.Range("AP32").Value = "HD76151"
$ws.Range("AP33").Value = "HD23356"
$ws.Range("AP34").Value = "HD20766"
$ws.Range("AP35").Value = "HD200779"
$ws.Range("AP36").Value = "HD222335"
$ws.Range("AP37").Value = "HD120467"
$ws.Range("AP38").Value = "HD181321"
$ws.Range("AP39").Value = "HD213042"
$ws.Range("AP40").Value = "HD173818"
$ws.Range("AP41").Value = "HD154088"
$ws.Range("AP42").Value = "HD192310"
$ws.Range("AP43").Value = "HD203040"
$ws.Range("AP44").Value = "HD14412"
$ws.Range("AP45").Value = "HD22049"
$ws.Range("AP46").Value = "HD177565"
$ws.Range("AP47").Value = "HD190248"
$ws.Range("AP49").Value = "HD85512"
$ws.Range("AP50").Value = "HD132683"
$ws.Range("AP51").Value = "HD69830"
$ws.Range("AP52").Value = "HD131977"
$ws.Range("AP53").Value = "HD120690"
$ws.Range("AP54").Value = "HD20280"
$ws.Range("AP55").Value = "HD209100"
$ws.Range("AP56").Value = "HD4391"
$ws.Range("AP57").Value = "HD21175"
$ws.Range("AP59").Value = "HD30876"
$ws.Range("AP60").Value = "HD158233"
$ws.Range("AP61").Value = "HD25004"
$ws.Range("AP62").Value = "HD16160"
$ws.Range("AP63").Value = "HD20807"
$ws.Range("AP64").Value = "HD1237A,HD1237"
$ws.Range("AP65").Value = "HD146233"
$ws.Range("AP66").Value = "HD41593"
$ws.Range("AP67").Value = "HD31560"
$ws.Range("AP68").Value = "HD111261B"
$ws.Range("AP69").Value = "HD139763"
$ws.Range("AP70").Value = "HD45088"
$ws.Range("AP71").Value = "HD222237"
$ws.Range("AP72").Value = "HD24916,HD24916A"
$ws.Range("AP73").Value = "HD32147"
$ws.Range("AP74").Value = "HD216803"
$ws.Range("AP75").Value = "HD65277,HD65277A"
$ws.Range("AP76").Value = "HD38858"
$ws.Range("AP77").Value = "HD30495"
$ws.Range("AP78").Value = "HD4628"
$ws.Range("AP79").Value = "HD166348"
$ws.Range("AP80").Value = "HD6101"
$ws.Range("AP81").Value = "HD221503"
$ws.Range("AP83").Value = "HD36003"
$ws.Range("AP85").Value = "HD75732"
$ws.Range("AP86").Value = "HD152391"
$ws.Range("AP87").Value = "HD274255"
$ws.Range("AP88").Value = "HD40307"
$ws.Range("AP89").Value = "HD53143"
$ws.Range("AP90").Value = "HD140538A,HD140538"
$ws.Range("AP91").Value = "HD72673"
$ws.Range("AP92").Value = "HD94765"
$ws.Range("AP93").Value = "HD17925"
$ws.Range("AP94").Value = "HD82558"
$ws.Range("AP95").Value = "HD145417"
$ws.Range("AP96").Value = "HD144628"
$ws.Range("AP97").Value = "HD109200"
$ws.Range("AP98").Value = "HD191391"
$ws.Range("AP99").Value = "HD125072"
$ws.Range("AP100").Value = "HD154577"
$ws.Range("AP101").Value = "HD168442"
$ws.Range("AP102").Value = "HD149661"
$ws.Range("AP103").Value = "HD35650"
$ws.Range("AP104").Value = "HD189733,HD189733A"
$ws.Range("AP105").Value = "HD192263"
$ws.Range("AP106").Value = "HD42807"
$ws.Range("AP107").Value = "HD211970"
$ws.Range("AP108").Value = "HD13445"
